$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 3758
$ws.Range("I3").Value = 3870
$ws.Range("H4").Value = 1666
$ws.Range("I4").Value = 914
$ws.Range("I5").Value = 360
$ws.Range("I6").Value = 4361
$ws.Range("H7").Value = 25976
$ws.Range("I7").Value = 13263

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I6").Value = 54
$ws.Range("I7").Value = 148

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I2").Value = 40
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I3").Value = 18
$ws.Range("I7").Value = 67

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 136
$ws.Range("I3").Value = 134
$ws.Range("I6").Value = 118
$ws.Range("I7").Value = 429

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I2").Value = 65
$ws.Range("I3").Value = 88
$ws.Range("I7").Value = 247

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 179
$ws.Range("I6").Value = 169
$ws.Range("I7").Value = 508

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 95
$ws.Range("I3").Value = 80
$ws.Range("I4").Value = 13
$ws.Range("I6").Value = 92
$ws.Range("I7").Value = 294

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 117
$ws.Range("I4").Value = 50
$ws.Range("I5").Value = 42
$ws.Range("I6").Value = 90
$ws.Range("I7").Value = 424
$ws.Range("H8").Value = 1696
$ws.Range("I8").Value = 801
$ws.Range("I14").Value = 67
$ws.Range("I18").Value = 90
$ws.Range("I19").Value = 357
$ws.Range("I20").Value = 324
$ws.Range("I23").Value = 126
$ws.Range("I25").Value = 64
$ws.Range("I27").Value = 126
$ws.Range("I29").Value = 855
$ws.Range("I33").Value = 599
$ws.Range("I36").Value = 183
$ws.Range("I37").Value = 429
$ws.Range("I42").Value = 461
$ws.Range("I43").Value = 116
$ws.Range("I44").Value = 95
$ws.Range("I45").Value = 26
$ws.Range("I47").Value = 93
$ws.Range("I50").Value = 57
$ws.Range("I51").Value = 127
$ws.Range("I53").Value = 144
$ws.Range("I54").Value = 303
$ws.Range("I55").Value = 145
$ws.Range("I63").Value = 49
$ws.Range("I64").Value = 118
$ws.Range("I65").Value = 294
$ws.Range("I67").Value = 508
$ws.Range("I68").Value = 42
$ws.Range("I76").Value = 201
$ws.Range("I78").Value = 192
$ws.Range("I79").Value = 357
$ws.Range("I85").Value = 598
$ws.Range("I86").Value = 83
$ws.Range("I87").Value = 23
$ws.Range("I88").Value = 119
$ws.Range("I89").Value = 148
$ws.Range("I90").Value = 167
$ws.Range("I91").Value = 163
$ws.Range("I93").Value = 75
$ws.Range("I96").Value = 146
$ws.Range("I97").Value = 100
$ws.Range("I98").Value = 87
$ws.Range("I99").Value = 247
$ws.Range("H101").Value = 25976
$ws.Range("I101").Value = 13263

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 217
$ws.Range("I6").Value = 189
$ws.Range("I7").Value = 599

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I3").Value = 58
$ws.Range("I6").Value = 153
$ws.Range("I7").Value = 303

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I3").Value = 290
$ws.Range("I5").Value = 35
$ws.Range("I6").Value = 237
$ws.Range("I7").Value = 855

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 135
$ws.Range("I3").Value = 102
$ws.Range("I4").Value = 16
$ws.Range("I6").Value = 96
$ws.Range("I7").Value = 357

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I6").Value = 28
$ws.Range("I7").Value = 95

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I4").Value = 24
$ws.Range("I7").Value = 201

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I6").Value = 147
$ws.Range("I7").Value = 598

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I3").Value = 27
$ws.Range("I7").Value = 90

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 126
$ws.Range("I3").Value = 157
$ws.Range("I6").Value = 124
$ws.Range("I7").Value = 461

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I2").Value = 40
$ws.Range("I6").Value = 78
$ws.Range("I7").Value = 192

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I2").Value = 45
$ws.Range("I6").Value = 50
$ws.Range("I7").Value = 145

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I3").Value = 43
$ws.Range("I7").Value = 126

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I3").Value = 57
$ws.Range("I7").Value = 163

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I3").Value = 113
$ws.Range("I6").Value = 106
$ws.Range("I7").Value = 357

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I2").Value = 33
$ws.Range("I7").Value = 118

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I2").Value = 91
$ws.Range("I7").Value = 324

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I5").Value = 2
$ws.Range("I7").Value = 90

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I3").Value = 56
$ws.Range("I6").Value = 55
$ws.Range("I7").Value = 183

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("I6").Value = 29
$ws.Range("I7").Value = 75

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("I2").Value = 20
$ws.Range("I7").Value = 64

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I4").Value = 7
$ws.Range("I7").Value = 93

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("I4").Value = 6
$ws.Range("I7").Value = 87

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I6").Value = 17
$ws.Range("I7").Value = 57

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I6").Value = 25
$ws.Range("I7").Value = 117

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I2").Value = 20
$ws.Range("I3").Value = 18
$ws.Range("I6").Value = 58
$ws.Range("I7").Value = 100

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I2").Value = 31
$ws.Range("I7").Value = 119

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 249
$ws.Range("H4").Value = 82
$ws.Range("H7").Value = 1696
$ws.Range("I7").Value = 801

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("I3").Value = 13
$ws.Range("I7").Value = 42

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I2").Value = 35
$ws.Range("I6").Value = 51
$ws.Range("I7").Value = 126

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I5").Value = 3
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 83

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I2").Value = 56
$ws.Range("I7").Value = 167

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I3").Value = 38
$ws.Range("I7").Value = 127

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("I4").Value = 6
$ws.Range("I7").Value = 42

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I6").Value = 68
$ws.Range("I7").Value = 116

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I2").Value = 30
$ws.Range("I7").Value = 144

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("I2").Value = 12
$ws.Range("I7").Value = 26

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I3").Value = 127
$ws.Range("I7").Value = 424

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("I3").Value = 14
$ws.Range("I7").Value = 50

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("I4").Value = 3
$ws.Range("I7").Value = 23
